$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row number, then a hashtable of Column -> new text value
# (G column always "17"; D/E updated per the refreshed price/volume snapshot)
$updates = @(
    @{ Row = 2; Cells = @{ "D"="282.22"; "E"="6.33%"; "G"="17" } }
    @{ Row = 3; Cells = @{ "D"="26.83"; "E"="0.64%"; "G"="17" } }
    @{ Row = 4; Cells = @{ "D"="4.946"; "E"="5.01%"; "G"="17" } }
    @{ Row = 5; Cells = @{ "D"="0.06401"; "E"="5.10%"; "G"="17" } }
    @{ Row = 6; Cells = @{ "D"="6.987"; "E"="3.61%"; "G"="17" } }
    @{ Row = 7; Cells = @{ "D"="3.354"; "E"="5.74%"; "G"="17" } }
    @{ Row = 8; Cells = @{ "D"="0.8863"; "E"="4.11%"; "G"="17" } }
    @{ Row = 9; Cells = @{ "D"="1.071"; "E"="17.81%"; "G"="17" } }
    @{ Row = 10; Cells = @{ "D"="0.1488"; "E"="5.96%"; "G"="17" } }
    @{ Row = 11; Cells = @{ "D"="0.05094"; "E"="1.05%"; "G"="17" } }
    @{ Row = 12; Cells = @{ "D"="0.07393"; "E"="4.11%"; "G"="17" } }
    @{ Row = 13; Cells = @{ "D"="0.03103"; "E"="-1.35%"; "G"="17" } }
    @{ Row = 14; Cells = @{ "D"="0.09054"; "E"="0.35%"; "G"="17" } }
    @{ Row = 15; Cells = @{ "D"="0.001564"; "E"="1.54%"; "G"="17" } }
    @{ Row = 16; Cells = @{ "D"="0.0006323"; "E"="4.64%"; "G"="17" } }
    @{ Row = 17; Cells = @{ "D"="0.005940"; "E"="-0.16%"; "G"="17" } }
    @{ Row = 18; Cells = @{ "D"="3.507"; "E"="1.67%"; "G"="17" } }
    @{ Row = 19; Cells = @{ "D"="2.299"; "E"="5.70%"; "G"="17" } }
    @{ Row = 20; Cells = @{ "G"="17" } }
    @{ Row = 21; Cells = @{ "D"="0.1291"; "E"="0.77%"; "G"="17" } }
    @{ Row = 22; Cells = @{ "D"="3.946"; "E"="-4.40%"; "G"="17" } }
    @{ Row = 23; Cells = @{ "D"="0.04342"; "E"="2.52%"; "G"="17" } }
    @{ Row = 24; Cells = @{ "D"="0.001176"; "E"="-0.30%"; "G"="17" } }
    @{ Row = 25; Cells = @{ "D"="0.003674"; "E"="-9.48%"; "G"="17" } }
    @{ Row = 26; Cells = @{ "D"="0.0001199"; "E"="-0.14%"; "G"="17" } }
    @{ Row = 27; Cells = @{ "D"="0.0001693"; "E"="0.65%"; "G"="17" } }
    @{ Row = 28; Cells = @{ "G"="17" } }
    @{ Row = 29; Cells = @{ "G"="17" } }
    @{ Row = 30; Cells = @{ "G"="17" } }
    @{ Row = 31; Cells = @{ "G"="17" } }
    @{ Row = 32; Cells = @{ "G"="17" } }
    @{ Row = 33; Cells = @{ "G"="17" } }
    @{ Row = 34; Cells = @{ "G"="17" } }
    @{ Row = 35; Cells = @{ "G"="17" } }
    @{ Row = 36; Cells = @{ "G"="17" } }
    @{ Row = 37; Cells = @{ "G"="17" } }
    @{ Row = 38; Cells = @{ "G"="17" } }
    @{ Row = 39; Cells = @{ "G"="17" } }
    @{ Row = 40; Cells = @{ "D"="0.04093"; "E"="4.19%"; "G"="17" } }
    @{ Row = 41; Cells = @{ "D"="0.006662"; "E"="58.69%"; "G"="17" } }
    @{ Row = 42; Cells = @{ "D"="0.1176"; "E"="5.73%"; "G"="17" } }
    @{ Row = 43; Cells = @{ "D"="0.002358"; "E"="11.96%"; "G"="17" } }
    @{ Row = 44; Cells = @{ "D"="0.01259"; "E"="9.48%"; "G"="17" } }
    @{ Row = 45; Cells = @{ "D"="0.00005266"; "E"="3.36%"; "G"="17" } }
    @{ Row = 46; Cells = @{ "E"="-0.07%"; "G"="17" } }
    @{ Row = 47; Cells = @{ "D"="2.355"; "E"="833.71%"; "G"="17" } }
    @{ Row = 48; Cells = @{ "D"="0.02248"; "E"="-8.14%"; "G"="17" } }
    @{ Row = 49; Cells = @{ "E"="-0.07%"; "G"="17" } }
    @{ Row = 50; Cells = @{ "E"="-0.14%"; "G"="17" } }
    @{ Row = 51; Cells = @{ "G"="17" } }
)

foreach ($update in $updates) {
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$($update.Row)"
        # Leading apostrophe forces text entry so Excel doesn't coerce
        # numeric-looking / percent-looking strings into numbers or dates.
        $ws.Range($cellRef).Value = "'$($update.Cells[$col])"
        $ws.Range($cellRef).Style = "Normal"
    }
}
